# Apply cryptos list update (data refresh) per commit:
# "Updated cryptos list on Sun May 21 09:10:56 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.011", "314.68").
# Force it to be stored as text so values/formatting exactly match the source data
# (trailing zeros, multi-dot grouping like "27.293.24", etc.) instead of being
# auto-coerced into Excel numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.293.24"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.831.47"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "314.68"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "0.4731"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").Value = "0.3686"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "0.07441"
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "20.47"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "1.888.49"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "0.07340"
$ws.Range("E13").Value = "  +3.16%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "94.20"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "5.429"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "6.560"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "0.000008791"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").Value = "27.644.14"
$ws.Range("E20").Value = "  +2.45%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").Value = "5.285"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "2.111.42"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").Value = "1.893"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "151.82"
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("D27").Value = "18.65"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "5.229"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "117.21"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "0.08991"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").Value = "0.7510"
$ws.Range("E32").Value = "  -1.17%  "
$ws.Range("D33").Value = "1.175"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").Value = "4.539"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "2.949"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "1.091"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "0.05346"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.430"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.967"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "7.249"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "0.5290"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "8.490"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "0.4931"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "10.55"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").Value = "105.19"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").Value = "1.011"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").Value = "1.668"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "0.06303"
$ws.Range("E51").Value = "  +0.09%  "

Write-Output "Applied 101 cell updates"
